$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update UserAccountId (column B) and UserDashboardId (column C) values
$ws.Range("B3").Value = 5
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 2
$ws.Range("B6").Value = 3
$ws.Range("B7").Value = 4

# Reflect the active cell selection change seen in the saved file
$ws.Range("B6").Select()
